$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force the Price (D) and Volume(1h) (E) columns to text format
# so values such as "3.60", "0.120" or "1.10" keep their exact literal
# representation instead of being coerced into numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# --- Cell value updates ---
$ws.Range("D2").Value = "41.942.18"
$ws.Range("E2").Value = "  -4.54%  "
$ws.Range("D3").Value = "2.221.89"
$ws.Range("E3").Value = "  -5.39%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "243.06"
$ws.Range("E5").Value = "  +0.93%  "
$ws.Range("E6").Value = "  -6.44%  "
$ws.Range("D7").Value = "67.81"
$ws.Range("E7").Value = "  -8.87%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -9.34%  "
$ws.Range("D10").Value = "0.0957"
$ws.Range("E10").Value = "  -5.25%  "
$ws.Range("D11").Value = "58.12"
$ws.Range("E11").Value = "  -2.55%  "
$ws.Range("D12").Value = "35.22"
$ws.Range("E12").Value = "  +5.80%  "
$ws.Range("E13").Value = "  -3.41%  "
$ws.Range("D14").Value = "6.67"
$ws.Range("E14").Value = "  -8.37%  "
$ws.Range("D15").Value = "2.553.78"
$ws.Range("E15").Value = "  -5.26%  "
$ws.Range("D16").Value = "14.71"
$ws.Range("E16").Value = "  -9.12%  "
$ws.Range("D17").Value = "0.844"
$ws.Range("E17").Value = "  -6.76%  "
$ws.Range("D18").Value = "2.231.02"
$ws.Range("D19").Value = "41.824.18"
$ws.Range("E19").Value = "  -4.58%  "
$ws.Range("D20").Value = "0.0₃0950"
$ws.Range("E20").Value = "  -8.15%  "
$ws.Range("D21").Value = "72.39"
$ws.Range("E21").Value = "  -7.55%  "
$ws.Range("E22").Value = "  -8.18%  "
$ws.Range("E23").Value = "  -7.54%  "
$ws.Range("D24").Value = "2.03"
$ws.Range("E24").Value = "  +10.25%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("D26").Value = "3.60"
$ws.Range("E26").Value = "  -5.07%  "
$ws.Range("D27").Value = "2.44"
$ws.Range("E27").Value = "  -2.52%  "
$ws.Range("D28").Value = "2.23"
$ws.Range("E28").Value = "  -3.15%  "
$ws.Range("D29").Value = "9.83"
$ws.Range("E29").Value = "  -5.95%  "
$ws.Range("D30").Value = "170.97"
$ws.Range("E30").Value = "  -2.91%  "
$ws.Range("E31").Value = "  -8.63%  "
$ws.Range("D32").Value = "0.120"
$ws.Range("E32").Value = "  -5.65%  "
$ws.Range("E33").Value = "  -6.99%  "
$ws.Range("E34").Value = "  -5.20%  "
$ws.Range("D35").Value = "5.17"
$ws.Range("E35").Value = "  -3.69%  "
$ws.Range("D36").Value = "4.64"
$ws.Range("E36").Value = "  -8.52%  "
$ws.Range("D37").Value = "3.85"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").Value = "22.75"
$ws.Range("E38").Value = "  +19.98%  "
$ws.Range("D39").Value = "0.0278"
$ws.Range("E39").Value = "  +2.17%  "
$ws.Range("D40").Value = "2.29"
$ws.Range("E40").Value = "  -3.94%  "
$ws.Range("D41").Value = "66.83"
$ws.Range("E41").Value = "  +3.26%  "
$ws.Range("E42").Value = "  -9.59%  "
$ws.Range("D43").Value = "4.94"
$ws.Range("D44").Value = "8.95"
$ws.Range("E44").Value = "  -2.27%  "
$ws.Range("E45").Value = "  -3.86%  "
$ws.Range("E46").Value = "  -4.48%  "
$ws.Range("E47").Value = "  +0.15%  "
$ws.Range("D48").Value = "4.56"
$ws.Range("E48").Value = "  +7.59%  "
$ws.Range("D49").Value = "1.17"
$ws.Range("E49").Value = "  -4.25%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D50").Value = "1.10"
$ws.Range("E50").Value = "  -4.32%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "2.81"
$ws.Range("E51").Value = "  -2.63%  "

# Restore the default (unstyled) appearance for the Price/Volume columns
# now that the text values have been written.
$ws.Range("D2:E51").Style = "Normal"
